# The commit swaps the presentation's applied theme from the custom
# "Integral" theme over to the stock "Office Theme" palette (the theme
# bound to the slide master, ppt/theme/theme1.xml). Re-create that by
# pushing the Office Theme's twelve scheme colours onto the slide
# master's ThemeColorScheme, in the same dk1/lt1/dk2/lt2/accent1-6/
# hlink/folHlink order the OOXML clrScheme uses.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

function Set-ThemeRGB {
    param(
        [int]$Index,
        [string]$Hex
    )
    $r = [Convert]::ToInt32($Hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4, 2), 16)
    $colorScheme.Item($Index).RGB = $r + ($g * 256) + ($b * 65536)
}

# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
Set-ThemeRGB 1 "000000"
Set-ThemeRGB 2 "FFFFFF"
Set-ThemeRGB 3 "44546A"
Set-ThemeRGB 4 "E7E6E6"
Set-ThemeRGB 5 "5B9BD5"
Set-ThemeRGB 6 "ED7D31"
Set-ThemeRGB 7 "A5A5A5"
Set-ThemeRGB 8 "FFC000"
Set-ThemeRGB 9 "4472C4"
Set-ThemeRGB 10 "70AD47"
Set-ThemeRGB 11 "0563C1"
Set-ThemeRGB 12 "954F72"
